# Generate Report for Handoff
# - Update status text from "Handed back: in sync with en-US" to "Ready for handoff"
#   on the Overview, zh-cn, and de-de sheets.
# - Update the related "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to reflect the new handoff generation time.
# - Narrow the (now shorter) status/date columns to match the new content width.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-12-16 09:27:55"

# Columns E (zh-cn) and F (de-de) shrink to match the new, shorter status text.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-12-16 09:27:41"

# Column C (Status) shrinks to match the new, shorter status text.
$zhcn.Columns.Item(3).ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-12-16 09:27:55"

# Column C (Status) shrinks to match the new, shorter status text.
$dede.Columns.Item(3).ColumnWidth = 16.33
